$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  "J2" = 7.753890582708475
  "L2" = 1.640365938249483
  "J3" = 8.041012805668172
  "L3" = 3.168949587529626
  "J4" = 8.143101297267171
  "L4" = 3.980459166710953
  "J5" = 7.654778896716499
  "L5" = 2.78096450573557
  "J6" = 7.355878924910716
  "L6" = 1.659683075975207
  "J7" = 7.260943120769692
  "L7" = 1.311750924950219
  "J8" = 7.093102683994399
  "L8" = 0.4963340482737552
  "J9" = 6.920808703097292
  "L9" = 0.3627216552763695
  "J10" = 6.983565512867058
  "L10" = 0.2270390024008371
  "J11" = 6.785013958627149
  "L11" = 1.595377886359134
  "J12" = 6.502521479593783
  "L12" = 4.649790550608018
  "J13" = 6.227900085676524
  "L13" = 7.869086316503365
  "J14" = 4.614765467253378
  "L14" = 0.9383607369115603
  "J15" = 4.547959936963774
  "L15" = 1.742011019337234
  "J16" = 4.65399111902201
  "L16" = 1.093777599097028
  "J17" = 4.738177902174329
  "L17" = 0.6336883705526255
  "J18" = 4.713101227274123
  "L18" = 0.9673456510876892
  "J19" = 4.604981200736534
  "L19" = 2.008066509097609
  "J20" = 4.574514404008129
  "L20" = 2.407420008542676
  "J21" = 4.57940770810222
  "L21" = 2.50500192761343
  "J22" = 4.616943023205891
  "L22" = 2.316605305096711
  "J23" = 4.653239657890016
  "L23" = 2.124627580255954
  "J24" = 4.494089830806232
  "L24" = 3.598599817182818
  "J25" = 4.479174050770306
  "L25" = 3.847416695453858
}

foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value = $values[$addr]
}
